$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item('展览')
$ws.Range('F2').Value = 22
$ws.Range('F4').Value = 5511
$ws.Range('F5').Value = 219
$ws.Range('B6').NumberFormat = "@"
$ws.Range('B6').Value = '2024-06-29'
$ws.Range('B6').Style = "Normal"
$ws.Range('C6').Value = '北京·Roll & Play跑团游戏大会'
$ws.Range('D6').Value = '中关村南大街12号 园艺创新中心'
$ws.Range('E6').Value = '2024.06.29 10:00-06.30 18:00'
$ws.Range('F6').Value = 6
$ws.Range('G6').Value = 78
$ws.Range('H6').Value = 'https://show.bilibili.com/platform/detail.html?id=86598'
$ws.Range('I6').Value = '//i0.hdslb.com/bfs/openplatform/202405/ZfmU4EqS1716911812871.jpeg'
$ws.Range('C7').Value = '北京·原神COSER全角色T台秀ONLY展'
$ws.Range('D7').Value = '广渠东路一号 创1958园区'
$ws.Range('E7').Value = '2024.06.29 09:30-06.30 18:30'
$ws.Range('F7').Value = 1252
$ws.Range('G7').Value = 69
$ws.Range('H7').Value = 'https://show.bilibili.com/platform/detail.html?id=83944'
$ws.Range('I7').Value = '//i1.hdslb.com/bfs/openplatform/202404/Lv3lJQKQ1714287707594.jpeg'
$ws.Range('C8').Value = '北京·日本经典动漫原稿＆吉卜力工作室原稿展'
$ws.Range('D8').Value = '酒仙桥路2号北京798艺术区红石广场东北角 北京第零空间'
$ws.Range('E8').Value = '2024.06.29 10:00-10.13 18:30'
$ws.Range('F8').Value = 9
$ws.Range('G8').Value = 29.9
$ws.Range('H8').Value = 'https://show.bilibili.com/platform/detail.html?id=87759'
$ws.Range('I8').Value = '//i1.hdslb.com/bfs/openplatform/202406/EmYCUgWo1718873950963.jpeg'
$ws.Range('B9').NumberFormat = "@"
$ws.Range('B9').Value = '2024-07-05'
$ws.Range('B9').Style = "Normal"
$ws.Range('C9').Value = '北京·樱桃小丸子限定快闪-人累拯救包'
$ws.Range('D9').Value = '西单大悦城 西单大悦城'
$ws.Range('E9').Value = '2024.07.05 10:00-08.04 22:00'
$ws.Range('F9').Value = 7
$ws.Range('G9').Value = 39
$ws.Range('H9').Value = 'https://show.bilibili.com/platform/detail.html?id=87724'
$ws.Range('I9').Value = '//i1.hdslb.com/bfs/openplatform/202406/5gq88LLf1718850914588.png'
$ws.Range('B10').NumberFormat = "@"
$ws.Range('B10').Value = '2024-07-06'
$ws.Range('B10').Style = "Normal"
$ws.Range('C10').Value = '北京·第六届璃樱动漫嘉年华'
$ws.Range('D10').Value = '永外高庄138号  大红门会展中心'
$ws.Range('E10').Value = '2024.07.06 10:00-07.06 17:00'
$ws.Range('F10').Value = 822
$ws.Range('G10').Value = 55
$ws.Range('H10').Value = 'https://show.bilibili.com/platform/detail.html?id=85472'
$ws.Range('I10').Value = '//i0.hdslb.com/bfs/openplatform/202405/5kwlHxZx1715063112027.png'
$ws.Range('B11').NumberFormat = "@"
$ws.Range('B11').Value = '2024-07-13'
$ws.Range('B11').Style = "Normal"
$ws.Range('C11').Value = '北京·Aw×SoReal二次元派对（Part2）'
$ws.Range('D11').Value = '石景山路68号首钢园内 首钢一高炉·SoReal科幻乐园'
$ws.Range('E11').Value = '2024.07.13 17:00-07.14 22:00'
$ws.Range('F11').Value = 26
$ws.Range('G11').Value = 108
$ws.Range('H11').Value = 'https://show.bilibili.com/platform/detail.html?id=87738'
$ws.Range('I11').Value = '//i1.hdslb.com/bfs/openplatform/202406/pQy6Nosb1718854298301.jpeg'
$ws.Range('C12').Value = '北京·Aw动漫游戏嘉年华8th-夏日奇幻之旅（Part1）'
$ws.Range('D12').Value = '石景山路68号 北京首钢会展中心'
$ws.Range('E12').Value = '2024.07.13 09:30-07.14 17:30'
$ws.Range('F12').Value = 6638
$ws.Range('G12').Value = 72
$ws.Range('H12').Value = 'https://show.bilibili.com/platform/detail.html?id=84800'
$ws.Range('I12').Value = '//i1.hdslb.com/bfs/openplatform/202405/hAKSdOQ91715586034060.jpeg'
$ws.Range('C13').Value = '北京·Aw动漫游戏嘉年华8th—coke老师专场见面会'
$ws.Range('E13').Value = '2024.07.13 10:00-07.13 17:00'
$ws.Range('F13').Value = 43
$ws.Range('G13').Value = 258
$ws.Range('H13').Value = 'https://show.bilibili.com/platform/detail.html?id=87337'
$ws.Range('I13').Value = '//i0.hdslb.com/bfs/openplatform/202406/ASxG2pZA1718245922365.png'
$ws.Range('C14').Value = '北京·Aw动漫游戏嘉年华8th—游马晃祐专场见面会'
$ws.Range('F14').Value = 90
$ws.Range('G14').Value = 588
$ws.Range('H14').Value = 'https://show.bilibili.com/platform/detail.html?id=87339'
$ws.Range('I14').Value = '//i1.hdslb.com/bfs/openplatform/202406/30OrK7QG1718248048340.png'
$ws.Range('C15').Value = '北京·Aw动漫游戏嘉年华8th—锦鲤专场见面会'
$ws.Range('F15').Value = 131
$ws.Range('G15').Value = 258
$ws.Range('H15').Value = 'https://show.bilibili.com/platform/detail.html?id=87342'
$ws.Range('I15').Value = '//i1.hdslb.com/bfs/openplatform/202406/aQHdZFWc1718250693994.png'
$ws.Range('C16').Value = '北京·GOJO超次元动漫游戏嘉年华14th'
$ws.Range('D16').Value = '小关路39号 北投购物公园'
$ws.Range('E16').Value = '2024.07.13 09:20-07.14 17:00'
$ws.Range('F16').Value = 6329
$ws.Range('G16').Value = 6.6
$ws.Range('H16').Value = 'https://show.bilibili.com/platform/detail.html?id=85225'
$ws.Range('I16').Value = '//i1.hdslb.com/bfs/openplatform/202406/mJt8McPp1718594709773.jpeg'
$ws.Range('B17').NumberFormat = "@"
$ws.Range('B17').Value = '2024-07-14'
$ws.Range('B17').Style = "Normal"
$ws.Range('C17').Value = '【大会员提前抢】北京·Aw动漫游戏嘉年华8th--谢莹签售礼包  '
$ws.Range('D17').Value = '石景山路68号 北京首钢会展中心'
$ws.Range('E17').Value = '2024.07.14 10:00-07.14 17:00'
$ws.Range('F17').Value = 115
$ws.Range('G17').Value = 59.9
$ws.Range('H17').Value = 'https://show.bilibili.com/platform/detail.html?id=87047'
$ws.Range('I17').Value = '//i2.hdslb.com/bfs/openplatform/202406/NZ43wLRW1717744995169.png'
$ws.Range('C18').Value = '北京·Aw动漫游戏嘉年华8th—帮我拍拍专场见面会'
$ws.Range('E18').Value = '2024.07.14 09:00-07.14 17:00'
$ws.Range('F18').Value = 265
$ws.Range('G18').Value = 398
$ws.Range('H18').Value = 'https://show.bilibili.com/platform/detail.html?id=86907'
$ws.Range('I18').Value = '//i2.hdslb.com/bfs/openplatform/202406/s2P9Isfw1717565356191.png'
$ws.Range('B19').NumberFormat = "@"
$ws.Range('B19').Value = '2024-07-20'
$ws.Range('B19').Style = "Normal"
$ws.Range('C19').Value = ' 北京·ICOS国际动漫节×CGF中国游戏节03'
$ws.Range('E19').Value = '2024.07.20 09:00-07.21 17:00'
$ws.Range('F19').Value = 4246
$ws.Range('G19').Value = 80
$ws.Range('H19').Value = 'https://show.bilibili.com/platform/detail.html?id=83931'
$ws.Range('I19').Value = '//i1.hdslb.com/bfs/openplatform/202404/sgFsCjWK1712558620744.jpeg'
$ws.Range('C20').Value = '【大会员提前抢】北京·ICOS内场-日本舞见鼻血姬'
$ws.Range('E20').Value = '2024.07.20 09:00-07.20 17:00'
$ws.Range('F20').Value = 5
$ws.Range('G20').Value = 168
$ws.Range('H20').Value = 'https://show.bilibili.com/platform/detail.html?id=86902'
$ws.Range('I20').Value = '//i1.hdslb.com/bfs/openplatform/202406/99SYO24h1717576009395.jpeg'
$ws.Range('C21').Value = '【大会员提前抢】北京·ICOS内场-谢安然'
$ws.Range('F21').Value = 48
$ws.Range('H21').Value = 'https://show.bilibili.com/platform/detail.html?id=86903'
$ws.Range('I21').Value = '//i0.hdslb.com/bfs/openplatform/202406/rLhhV7bQ1717576183936.jpeg'
$ws.Range('F22').Value = 4187
$ws.Range('F24').Value = 218
$ws.Range('F26').Value = 293
$ws.Range('F28').Value = 223
$ws.Range('F33').Value = 7604
$ws.Range('F34').Value = 48
$ws.Range('F35').Value = 1280
$ws.Range('F36').Value = 633
$ws.Range('F37').Value = 114
$ws.Range('F38').Value = 980
$ws.Range('F39').Value = 66
$ws.Range('F40').Value = 1517
$ws.Range('F41').Value = 201
$ws.Range('F42').Value = 853
$ws.Range('F44').Value = 3752
$ws.Range('F45').Value = 336
$ws.Range('F47').Value = 103
$ws.Range('F49').Value = 1050

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item('演出')
$ws.Range('F13').Value = 153
$ws.Range('F18').Value = 71
$ws.Range('F21').Value = 859

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F3').Value = 22
$ws.Range('F5').Value = 5511
$ws.Range('F6').Value = 219
$ws.Range('C7').Value = '北京·仲夏绮罗日 Anisong Live Party'
$ws.Range('D7').Value = '大江胡同121号2幢负1层 北京门空间 TheDoorLiveHouse'
$ws.Range('E7').Value = '2024.06.23 13:00-06.23 16:30'
$ws.Range('F7').Value = 101
$ws.Range('G7').Value = 98
$ws.Range('H7').Value = 'https://show.bilibili.com/platform/detail.html?id=85364'
$ws.Range('I7').Value = '//i0.hdslb.com/bfs/openplatform/202405/PFM2Be6V1715240437688.jpeg'
$ws.Range('B8').NumberFormat = "@"
$ws.Range('B8').Value = '2024-06-28'
$ws.Range('B8').Style = "Normal"
$ws.Range('C8').Value = '北京·《国风大赏》大型国潮音乐会×郑州歌舞剧院《唐宫夜宴》'
$ws.Range('D8').Value = '西直门外大街135号（北京展览馆内） 北京展览馆剧场'
$ws.Range('E8').Value = '2024.06.28 19:30-06.28 21:00'
$ws.Range('F8').Value = 69
$ws.Range('G8').Value = 162
$ws.Range('H8').Value = 'https://show.bilibili.com/platform/detail.html?id=82587'
$ws.Range('I8').Value = '//i2.hdslb.com/bfs/openplatform/202403/VZcJ2SJ51709882503997.jpeg'
$ws.Range('F10').Value = 1252
$ws.Range('F14').Value = 822
$ws.Range('F15').Value = 26
$ws.Range('F16').Value = 6638
$ws.Range('F18').Value = 90
$ws.Range('F19').Value = 131
$ws.Range('F20').Value = 6329
$ws.Range('F22').Value = 265
$ws.Range('F23').Value = 4246
$ws.Range('F24').Value = 4187
$ws.Range('F26').Value = 218
$ws.Range('F27').Value = 293
$ws.Range('F29').Value = 223
$ws.Range('F31').Value = 153
$ws.Range('F32').Value = 7604
$ws.Range('F33').Value = 48
$ws.Range('F34').Value = 1280
$ws.Range('F35').Value = 633
$ws.Range('F36').Value = 114
$ws.Range('F37').Value = 980
$ws.Range('F38').Value = 66
$ws.Range('F39').Value = 1517
$ws.Range('F40').Value = 201
$ws.Range('F41').Value = 853
$ws.Range('F43').Value = 3752
$ws.Range('F44').Value = 336
$ws.Range('F46').Value = 103
$ws.Range('F48').Value = 1050
